$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 317.1111
$ws.Range("I33").Value = 120.15385
$ws.Range("K33").Value = 120.15385
$ws.Range("M33").Value = 108.84615
$ws.Range("H62").Value = 6351.8667
$ws.Range("I62").Value = 4921.5
$ws.Range("J62").Value = 7305.4443
$ws.Range("K62").Value = 4921.5
$ws.Range("L62").Value = 7305.4443
$ws.Range("M62").Value = -4297.5
$ws.Range("N62").Value = -8553.444299999999
$ws.Range("H65").Value = 6351.8667
$ws.Range("I65").Value = 4921.5
$ws.Range("J65").Value = 7305.4443
$ws.Range("K65").Value = 24607.5
$ws.Range("L65").Value = 36527.2215
$ws.Range("M65").Value = -21487.5
$ws.Range("N65").Value = -42767.2215
$ws.Range("H107").Value = 429.3684
$ws.Range("I107").Value = 398.33334
$ws.Range("K107").Value = 398.33334
$ws.Range("M107").Value = 1521.66666
$ws.Range("H115").Value = 9634.799999999999
$ws.Range("I115").Value = 9634.799999999999
$ws.Range("K115").Value = 28904.4
$ws.Range("M115").Value = -27337.4
$ws.Range("H116").Value = 5684
$ws.Range("I116").Value = 3696.6
$ws.Range("K116").Value = 3696.6
$ws.Range("M116").Value = -254.5999999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2308.825
$ws.Range("I132").Value = 2270.7715
$ws.Range("K132").Value = 6812.314499999999
$ws.Range("M132").Value = -4282.314499999999
$ws.Range("H138").Value = 5917.559
$ws.Range("J138").Value = 5467.419
$ws.Range("L138").Value = 16402.257
$ws.Range("N138").Value = -26682.257
$ws.Range("H32").Value = 24383.531
$ws.Range("I32").Value = 14390
$ws.Range("K32").Value = 14390
$ws.Range("M32").Value = -14103
$ws.Range("H88").Value = 4200
$ws.Range("I88").Value = 1333.3334
$ws.Range("J88").Value = 7066.6665
$ws.Range("K88").Value = 1333.3334
$ws.Range("L88").Value = 7066.6665
$ws.Range("M88").Value = -927.3334
$ws.Range("N88").Value = -7878.6665
$ws.Range("H91").Value = 4200
$ws.Range("I91").Value = 1333.3334
$ws.Range("J91").Value = 7066.6665
$ws.Range("K91").Value = 1333.3334
$ws.Range("L91").Value = 7066.6665
$ws.Range("M91").Value = 70.66660000000002
$ws.Range("N91").Value = -9874.666499999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 5806.5
$ws.Range("I86").Value = 4806
$ws.Range("J86").Value = 6807
$ws.Range("K86").Value = 4806
$ws.Range("L86").Value = 6807
$ws.Range("M86").Value = -3683
$ws.Range("N86").Value = -9053
$ws.Range("H89").Value = 5806.5
$ws.Range("I89").Value = 4806
$ws.Range("J89").Value = 6807
$ws.Range("K89").Value = 24030
$ws.Range("L89").Value = 34035
$ws.Range("M89").Value = -18414
$ws.Range("N89").Value = -45267
$ws.Range("H94").Value = 1474.7273
$ws.Range("I94").Value = 1524.6666
$ws.Range("J94").Value = 1250
$ws.Range("K94").Value = 1524.6666
$ws.Range("L94").Value = 1250
$ws.Range("M94").Value = -1073.6666
$ws.Range("N94").Value = -2152

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4065.628
$ws.Range("I31").Value = 2937.0833
$ws.Range("J31").Value = 5491.1577
$ws.Range("K31").Value = 2937.0833
$ws.Range("L31").Value = 5491.1577
$ws.Range("M31").Value = -2642.0833
$ws.Range("N31").Value = -6081.1577
$ws.Range("H34").Value = 4065.628
$ws.Range("I34").Value = 2937.0833
$ws.Range("J34").Value = 5491.1577
$ws.Range("K34").Value = 2937.0833
$ws.Range("L34").Value = 5491.1577
$ws.Range("M34").Value = -2735.0833
$ws.Range("N34").Value = -5895.1577
$ws.Range("H99").Value = 10829.223
$ws.Range("I99").Value = 7244.5713
$ws.Range("J99").Value = 14689.615
$ws.Range("K99").Value = 7244.5713
$ws.Range("L99").Value = 14689.615
$ws.Range("M99").Value = -5746.5713
$ws.Range("N99").Value = -17685.615
$ws.Range("H107").Value = 932.3889
$ws.Range("I107").Value = 290.8
$ws.Range("K107").Value = 290.8
$ws.Range("M107").Value = 1629.2
$ws.Range("H126").Value = 10829.223
$ws.Range("I126").Value = 7244.5713
$ws.Range("J126").Value = 14689.615
$ws.Range("K126").Value = 21733.7139
$ws.Range("L126").Value = 44068.845
$ws.Range("M126").Value = -19263.7139
$ws.Range("N126").Value = -49008.845

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H59").Value = 8218.111000000001
$ws.Range("I59").Value = 7540.75
$ws.Range("K59").Value = 22622.25
$ws.Range("M59").Value = -22082.25
$ws.Range("H60").Value = 915.55554
$ws.Range("J60").Value = 832.8333
$ws.Range("L60").Value = 2498.4999
$ws.Range("N60").Value = -3000.4999
$ws.Range("H70").Value = 1937.1428
$ws.Range("I70").Value = 1426.6666
$ws.Range("K70").Value = 4279.9998
$ws.Range("M70").Value = -3964.9998
$ws.Range("H73").Value = 1937.1428
$ws.Range("I73").Value = 1426.6666
$ws.Range("K73").Value = 4279.9998
$ws.Range("M73").Value = -3187.9998
$ws.Range("H81").Value = 4007
$ws.Range("J81").Value = 4007
$ws.Range("L81").Value = 12021
$ws.Range("N81").Value = -14267
$ws.Range("H84").Value = 4007
$ws.Range("J84").Value = 4007
$ws.Range("L84").Value = 36063
$ws.Range("N84").Value = -47295
$ws.Range("H121").Value = 957.1667
$ws.Range("I121").Value = 585.6
$ws.Range("J121").Value = 1222.5714
$ws.Range("K121").Value = 1756.8
$ws.Range("L121").Value = 3667.7142
$ws.Range("M121").Value = -446.8000000000002
$ws.Range("N121").Value = -6287.7142

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 3697.8572
$ws.Range("I82").Value = 3777.2
$ws.Range("J82").Value = 3499.5
$ws.Range("K82").Value = 3777.2
$ws.Range("L82").Value = 3499.5
$ws.Range("M82").Value = -3416.2
$ws.Range("N82").Value = -4221.5
$ws.Range("H85").Value = 3697.8572
$ws.Range("I85").Value = 3777.2
$ws.Range("J85").Value = 3499.5
$ws.Range("K85").Value = 3777.2
$ws.Range("L85").Value = 3499.5
$ws.Range("M85").Value = -2529.2
$ws.Range("N85").Value = -5995.5
$ws.Range("H122").Value = 8299.444
$ws.Range("I122").Value = 6939
$ws.Range("K122").Value = 20817
$ws.Range("M122").Value = -18367

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 5080
$ws.Range("I58").Value = 5080
$ws.Range("K58").Value = 5080
$ws.Range("M58").Value = -4772
$ws.Range("H62").Value = 6566.5
$ws.Range("I62").Value = 4799
$ws.Range("K62").Value = 4799
$ws.Range("M62").Value = -4175
$ws.Range("H65").Value = 6566.5
$ws.Range("I65").Value = 4799
$ws.Range("K65").Value = 23995
$ws.Range("M65").Value = -20875
$ws.Range("H126").Value = 97952.09
$ws.Range("I126").Value = 172913.83
$ws.Range("J126").Value = 7998
$ws.Range("K126").Value = 518741.49
$ws.Range("L126").Value = 23994
$ws.Range("M126").Value = -516271.49
$ws.Range("N126").Value = -28934
$ws.Range("H135").Value = 80000
$ws.Range("J135").Value = 80000
$ws.Range("L135").Value = 80000
$ws.Range("N135").Value = -90140
$ws.Range("H136").Value = 254847.25
$ws.Range("I136").Value = 1400
$ws.Range("K136").Value = 4200
$ws.Range("M136").Value = -1650
$ws.Range("H141").Value = 74300
$ws.Range("J141").Value = 74300
$ws.Range("L141").Value = 74300
$ws.Range("N141").Value = -84660
